$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20
$data[0,0] = "ECs"
$data[0,1] = "Amelx"
$data[0,2] = "Lamp2"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.04979866666666666
$data[0,7] = 0.149396
$data[0,8] = 0.04800668898255549
$data[0,9] = 0.04800668898255549
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 47.13631833333333
$data[0,13] = 141.408955
$data[0,14] = 0.05713678551147761
$data[0,15] = 0.05843911332169659
$data[0,16] = 2.347325804575555
$data[0,17] = 21.12593224118
$data[0,18] = 0.002742947891512489
$data[0,19] = 0.002805468337651004
$data[1,0] = "ECs"
$data[1,1] = "Amelx"
$data[1,2] = "Lamp2"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.04979866666666666
$data[1,7] = 0.149396
$data[1,8] = 0.04800668898255549
$data[1,9] = 0.04800668898255549
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 192.2656683333333
$data[1,13] = 576.797005
$data[1,14] = 0.2330568580918209
$data[1,15] = 0.2383689600054692
$data[1,16] = 9.574573928775555
$data[1,17] = 86.17116535898001
$data[1,18] = 0.01118828810166562
$data[1,19] = 0.01144330452607777
$data[2,0] = "ECs"
$data[2,1] = "Amelx"
$data[2,2] = "Lamp2"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.04979866666666666
$data[2,7] = 0.149396
$data[2,8] = 0.04800668898255549
$data[2,9] = 0.04800668898255549
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 277.4609476666666
$data[2,13] = 832.3828429999999
$data[2,14] = 0.3363272146656126
$data[2,15] = 0.3439931741883884
$data[2,16] = 13.81718524586978
$data[2,17] = 124.354667212828
$data[2,18] = 0.01614595599082124
$data[2,19] = 0.016513973325384
$data[3,0] = "ECs"
$data[3,1] = "Amelx"
$data[3,2] = "Lamp2"
$data[3,3] = "MuSCs"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.04979866666666666
$data[3,7] = 0.149396
$data[3,8] = 0.04800668898255549
$data[3,9] = 0.04800668898255549
$data[3,10] = 2
$data[3,11] = 1
$data[3,12] = 55.1541005
$data[3,13] = 110.308201
$data[3,14] = 0.06685562474484689
$data[3,15] = 0.04558631706564471
$data[3,16] = 2.746600666099333
$data[3,17] = 16.479603996596
$data[3,18] = 0.003209517183860306
$data[3,19] = 0.002188448145230568
$data[4,0] = "ECs"
$data[4,1] = "Amelx"
$data[4,2] = "Lamp2"
$data[4,3] = "Resolving-Mac"
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.04979866666666666
$data[4,7] = 0.149396
$data[4,8] = 0.04800668898255549
$data[4,9] = 0.04800668898255549
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 252.9561923333333
$data[4,13] = 758.8685770000001
$data[4,14] = 0.3066235169862421
$data[4,15] = 0.313612435418801
$data[4,16] = 12.59688110327689
$data[4,17] = 113.371929929492
$data[4,18] = 0.01471997981469584
$data[4,19] = 0.01505549464821215
$data[5,0] = "FAPs"
$data[5,1] = "Amelx"
$data[5,2] = "Lamp2"
$data[5,3] = "ECs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.809644
$data[5,7] = 2.428932
$data[5,8] = 0.7805094050963647
$data[5,9] = 0.7805094050963647
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 47.13631833333333
$data[5,13] = 141.408955
$data[5,14] = 0.05713678551147761
$data[5,15] = 0.05843911332169659
$data[5,16] = 38.16363732067333
$data[5,17] = 343.47273588606
$data[5,18] = 0.04459579846868198
$data[5,19] = 0.04561227757307644
$data[6,0] = "FAPs"
$data[6,1] = "Amelx"
$data[6,2] = "Lamp2"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.809644
$data[6,7] = 2.428932
$data[6,8] = 0.7805094050963647
$data[6,9] = 0.7805094050963647
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 192.2656683333333
$data[6,13] = 576.797005
$data[6,14] = 0.2330568580918209
$data[6,15] = 0.2383689600054692
$data[6,16] = 155.6667447720733
$data[6,17] = 1401.00070294866
$data[6,18] = 0.181903069662875
$data[6,19] = 0.1860492151673079
$data[7,0] = "FAPs"
$data[7,1] = "Amelx"
$data[7,2] = "Lamp2"
$data[7,3] = "Inflammatory-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.809644
$data[7,7] = 2.428932
$data[7,8] = 0.7805094050963647
$data[7,9] = 0.7805094050963647
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 277.4609476666666
$data[7,13] = 832.3828429999999
$data[7,14] = 0.3363272146656126
$data[7,15] = 0.3439931741883884
$data[7,16] = 224.6445915126307
$data[7,17] = 2021.801323613676
$data[7,18] = 0.2625065542363746
$data[7,19] = 0.2684899077429891
$data[8,0] = "FAPs"
$data[8,1] = "Amelx"
$data[8,2] = "Lamp2"
$data[8,3] = "MuSCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.809644
$data[8,7] = 2.428932
$data[8,8] = 0.7805094050963647
$data[8,9] = 0.7805094050963647
$data[8,10] = 2
$data[8,11] = 1
$data[8,12] = 55.1541005
$data[8,13] = 110.308201
$data[8,14] = 0.06685562474484689
$data[8,15] = 0.04558631706564471
$data[8,16] = 44.655186545222
$data[8,17] = 267.931119271332
$data[8,18] = 0.05218144389694625
$data[8,19] = 0.03558054921344062
$data[9,0] = "FAPs"
$data[9,1] = "Amelx"
$data[9,2] = "Lamp2"
$data[9,3] = "Resolving-Mac"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.809644
$data[9,7] = 2.428932
$data[9,8] = 0.7805094050963647
$data[9,9] = 0.7805094050963647
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 252.9561923333333
$data[9,13] = 758.8685770000001
$data[9,14] = 0.3066235169862421
$data[9,15] = 0.313612435418801
$data[9,16] = 204.8044633855293
$data[9,17] = 1843.240170469764
$data[9,18] = 0.2393225388314869
$data[9,19] = 0.2447774553995505
$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "Amelx"
$data[10,2] = "Lamp2"
$data[10,3] = "ECs"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.036349
$data[10,7] = 0.109047
$data[10,8] = 0.0350410011879885
$data[10,9] = 0.0350410011879885
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 47.13631833333333
$data[10,13] = 141.408955
$data[10,14] = 0.05713678551147761
$data[10,15] = 0.05843911332169659
$data[10,16] = 1.713358035098333
$data[10,17] = 15.420222315885
$data[10,18] = 0.002002130168985531
$data[10,19] = 0.002047765039330564
$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "Amelx"
$data[11,2] = "Lamp2"
$data[11,3] = "FAPs"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.036349
$data[11,7] = 0.109047
$data[11,8] = 0.0350410011879885
$data[11,9] = 0.0350410011879885
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 192.2656683333333
$data[11,13] = 576.797005
$data[11,14] = 0.2330568580918209
$data[11,15] = 0.2383689600054692
$data[11,16] = 6.988664778248333
$data[11,17] = 62.897983004235
$data[11,18] = 0.00816654564126436
$data[11,19] = 0.008352687010731229
$data[12,0] = "Inflammatory-Mac"
$data[12,1] = "Amelx"
$data[12,2] = "Lamp2"
$data[12,3] = "Inflammatory-Mac"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.036349
$data[12,7] = 0.109047
$data[12,8] = 0.0350410011879885
$data[12,9] = 0.0350410011879885
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 277.4609476666666
$data[12,13] = 832.3828429999999
$data[12,14] = 0.3363272146656126
$data[12,15] = 0.3439931741883884
$data[12,16] = 10.08542798673566
$data[12,17] = 90.76885188062099
$data[12,18] = 0.01178524232865059
$data[12,19] = 0.01205386522539525
$data[13,0] = "Inflammatory-Mac"
$data[13,1] = "Amelx"
$data[13,2] = "Lamp2"
$data[13,3] = "MuSCs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.036349
$data[13,7] = 0.109047
$data[13,8] = 0.0350410011879885
$data[13,9] = 0.0350410011879885
$data[13,10] = 2
$data[13,11] = 1
$data[13,12] = 55.1541005
$data[13,13] = 110.308201
$data[13,14] = 0.06685562474484689
$data[13,15] = 0.04558631706564471
$data[13,16] = 2.0047963990745
$data[13,17] = 12.028778394447
$data[13,18] = 0.002342688026107893
$data[13,19] = 0.001597390190453277
$data[14,0] = "Inflammatory-Mac"
$data[14,1] = "Amelx"
$data[14,2] = "Lamp2"
$data[14,3] = "Resolving-Mac"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.036349
$data[14,7] = 0.109047
$data[14,8] = 0.0350410011879885
$data[14,9] = 0.0350410011879885
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 252.9561923333333
$data[14,13] = 758.8685770000001
$data[14,14] = 0.3066235169862421
$data[14,15] = 0.313612435418801
$data[14,16] = 9.194704635124333
$data[14,17] = 82.75234171611901
$data[14,18] = 0.01074439502298012
$data[14,19] = 0.01098929372207817
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Amelx"
$data[15,2] = "Lamp2"
$data[15,3] = "ECs"
$data[15,4] = 1
$data[15,5] = 0.3333333333333333
$data[15,6] = 0.141536
$data[15,7] = 0.424608
$data[15,8] = 0.1364429047330914
$data[15,9] = 0.1364429047330914
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 47.13631833333333
$data[15,13] = 141.408955
$data[15,14] = 0.05713678551147761
$data[15,15] = 0.05843911332169659
$data[15,16] = 6.671485951626666
$data[15,17] = 60.04337356463999
$data[15,18] = 0.007795908982297616
$data[15,19] = 0.007973602371638579
$data[16,0] = "Resolving-Mac"
$data[16,1] = "Amelx"
$data[16,2] = "Lamp2"
$data[16,3] = "FAPs"
$data[16,4] = 1
$data[16,5] = 0.3333333333333333
$data[16,6] = 0.141536
$data[16,7] = 0.424608
$data[16,8] = 0.1364429047330914
$data[16,9] = 0.1364429047330914
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 192.2656683333333
$data[16,13] = 576.797005
$data[16,14] = 0.2330568580918209
$data[16,15] = 0.2383689600054692
$data[16,16] = 27.21251363322667
$data[16,17] = 244.91262269904
$data[16,18] = 0.03179895468601592
$data[16,19] = 0.03252375330135231
$data[17,0] = "Resolving-Mac"
$data[17,1] = "Amelx"
$data[17,2] = "Lamp2"
$data[17,3] = "Inflammatory-Mac"
$data[17,4] = 1
$data[17,5] = 0.3333333333333333
$data[17,6] = 0.141536
$data[17,7] = 0.424608
$data[17,8] = 0.1364429047330914
$data[17,9] = 0.1364429047330914
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 277.4609476666666
$data[17,13] = 832.3828429999999
$data[17,14] = 0.3363272146656126
$data[17,15] = 0.3439931741883884
$data[17,16] = 39.27071268894933
$data[17,17] = 353.4364142005439
$data[17,18] = 0.04588946210976615
$data[17,19] = 0.04693542789461999
$data[18,0] = "Resolving-Mac"
$data[18,1] = "Amelx"
$data[18,2] = "Lamp2"
$data[18,3] = "MuSCs"
$data[18,4] = 1
$data[18,5] = 0.3333333333333333
$data[18,6] = 0.141536
$data[18,7] = 0.424608
$data[18,8] = 0.1364429047330914
$data[18,9] = 0.1364429047330914
$data[18,10] = 2
$data[18,11] = 1
$data[18,12] = 55.1541005
$data[18,13] = 110.308201
$data[18,14] = 0.06685562474484689
$data[18,15] = 0.04558631706564471
$data[18,16] = 7.806290768367999
$data[18,17] = 46.837744610208
$data[18,18] = 0.009121975637932453
$data[18,19] = 0.00621992951652026
$data[19,0] = "Resolving-Mac"
$data[19,1] = "Amelx"
$data[19,2] = "Lamp2"
$data[19,3] = "Resolving-Mac"
$data[19,4] = 1
$data[19,5] = 0.3333333333333333
$data[19,6] = 0.141536
$data[19,7] = 0.424608
$data[19,8] = 0.1364429047330914
$data[19,9] = 0.1364429047330914
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 252.9561923333333
$data[19,13] = 758.8685770000001
$data[19,14] = 0.3066235169862421
$data[19,15] = 0.313612435418801
$data[19,16] = 35.80240763809067
$data[19,17] = 322.221668742816
$data[19,18] = 0.04183660331707926
$data[19,19] = 0.04279019164896024

$ws.Range("A2:T21").Value = $data
